$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each crypto row per the latest data refresh.
# D-column numeric-looking values are forced to remain text (matching the original inlineStr/shared-string
# cell type) by briefly applying a Text number format, then resetting the cell style back to Normal so no
# stray style index is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.337.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.060.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.03%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.96"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.62%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.08%  "

$ws.Range("E11").Value = "  +2.20%  "

$ws.Range("E12").Value = "  +3.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.363.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.776"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.31%  "

$ws.Range("E17").Value = "  +2.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.063.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.538.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +16.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0815"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.23%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "

$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("E32").Value = "  +1.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("E34").Value = "  +2.82%  "

$ws.Range("E35").Value = "  +11.33%  "

$ws.Range("E36").Value = "  +5.75%  "

$ws.Range("E37").Value = "  +6.28%  "

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +30.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0984"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.476.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.98%  "

$ws.Range("E46").Value = "  +6.85%  "

$ws.Range("E47").Value = "  +4.59%  "

$ws.Range("E48").Value = "  +7.34%  "

$ws.Range("E49").Value = "  +3.32%  "

$ws.Range("E50").Value = "  +6.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.17%  "
